# Swap the "HẠNG 2" (Hà Nội) and "HẠNG 3" (QDO) team blocks in the
# "KẾT QUẢ ĐỒNG ĐỘI" standings sheet, along with their score data,
# so that QDO becomes rank 2 and Hà Nội becomes rank 3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Section header row 6: was "HẠNG 2 ĐỘI: Hà Nội" -> "HẠNG 2 ĐỘI: QDO"
$ws.Range("A6").Value = "HẠNG 2 ĐỘI: QDO"

# --- Row 7 (was Hà Nội 1 / Hà Nội): becomes QDO 4 / QDO
$ws.Range("B7").Value = "QDO 4"
$ws.Range("C7").Value = "QDO"
$ws.Range("A7").Value = 12
$ws.Range("D7").Value = 3
$ws.Range("E7").Value = 8
$ws.Range("G7").Value = 12

# --- Row 8 (was Hà Nội 2 / Hà Nội): becomes QDO 2 / QDO
$ws.Range("B8").Value = "QDO 2"
$ws.Range("C8").Value = "QDO"
$ws.Range("A8").Value = 10
$ws.Range("D8").Value = 7
$ws.Range("E8").Value = 4

# --- Section header row 9: was "HẠNG 3 ĐỘI: QDO" -> "HẠNG 3 ĐỘI: Hà Nội"
$ws.Range("A9").Value = "HẠNG 3 ĐỘI: Hà Nội"

# --- Row 10 (was QDO 4 / QDO): becomes Hà Nội 1 / Hà Nội
$ws.Range("B10").Value = "Hà Nội 1"
$ws.Range("C10").Value = "Hà Nội"
$ws.Range("A10").Value = 1
$ws.Range("D10").Value = 2
$ws.Range("E10").Value = 6
$ws.Range("G10").Value = 9

# --- Row 11 (was QDO 2 / QDO): becomes Hà Nội 2 / Hà Nội
$ws.Range("B11").Value = "Hà Nội 2"
$ws.Range("C11").Value = "Hà Nội"
$ws.Range("A11").Value = 2
$ws.Range("D11").Value = 8
$ws.Range("E11").Value = 3
